{"js": "// Replace every arithmetic-problem cell in the single table with its\n// updated \"after\" value, preserving table/row/cell structure and the\n// existing run formatting (font, size) already present in each cell.\nconst newValues = [\n  [\"19+73=\", \"37+55=\", \"85-27=\", \"16+77=\", \"94-88=\"],\n  [\"86-47=\", \"3+38=\", \"65+27=\", \"60-26=\", \"52-38=\"],\n  [\"92-46=\", \"4+28=\", \"83-66=\", \"92-3=\", \"81-29=\"],\n  [\"49+25=\", \"53-14=\", \"28+13=\", \"48+34=\", \"4+78=\"],\n  [\"4+19=\", \"17+34=\", \"55+27=\", \"27+28=\", \"49+33=\"],\n  [\"94-5=\", \"71-12=\", \"25+66=\", \"33+8=\", \"71-5=\"],\n  [\"46+15=\", \"20-9=\", \"60-31=\", \"35-28=\", \"80-47=\"],\n  [\"91-45=\", \"69+19=\", \"5+67=\", \"23+9=\", \"9+13=\"],\n  [\"27+46=\", \"7+45=\", \"74-57=\", \"32-6=\", \"88-59=\"],\n  [\"58+28=\", \"12+69=\", \"49+42=\", \"42-39=\", \"37+46=\"],\n  [\"54+37=\", \"36-7=\", \"27+25=\", \"75-68=\", \"8+84=\"],\n  [\"56-28=\", \"58+34=\", \"49+35=\", \"12-7=\", \"93-69=\"],\n  [\"76+8=\", \"71-67=\", \"27-19=\", \"29+36=\", \"74-28=\"],\n  [\"74-6=\", \"8+56=\", \"29+35=\", \"38+15=\", \"68-49=\"],\n  [\"44-18=\", \"98-79=\", \"36+15=\", \"70-33=\", \"80-48=\"],\n  [\"44-19=\", \"34+37=\", \"90-65=\", \"16+58=\", \"73-15=\"],\n  [\"86+8=\", \"83-66=\", \"25+57=\", \"37+7=\", \"9+84=\"],\n  [\"16+19=\", \"23-6=\", \"17+69=\", \"15+78=\", \"42-26=\"],\n  [\"46+15=\", \"55+27=\", \"79+3=\", \"52-38=\", \"71-22=\"],\n  [\"55+39=\", \"53-35=\", \"97-29=\", \"13+58=\", \"4+17=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Assigning the full 2-D array of strings rewrites each cell's text\n// content in place (keeping each cell's existing paragraph/run\n// formatting), which matches how the diff only touches the <w:t> text.\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace every arithmetic-problem cell in the single table with its\n# updated \"after\" value, preserving table/row/cell structure and the\n# existing run formatting (font, size) already present in each cell.\n\n$newValues = @(\n    @(\"19+73=\", \"37+55=\", \"85-27=\", \"16+77=\", \"94-88=\"),\n    @(\"86-47=\", \"3+38=\", \"65+27=\", \"60-26=\", \"52-38=\"),\n    @(\"92-46=\", \"4+28=\", \"83-66=\", \"92-3=\", \"81-29=\"),\n    @(\"49+25=\", \"53-14=\", \"28+13=\", \"48+34=\", \"4+78=\"),\n    @(\"4+19=\", \"17+34=\", \"55+27=\", \"27+28=\", \"49+33=\"),\n    @(\"94-5=\", \"71-12=\", \"25+66=\", \"33+8=\", \"71-5=\"),\n    @(\"46+15=\", \"20-9=\", \"60-31=\", \"35-28=\", \"80-47=\"),\n    @(\"91-45=\", \"69+19=\", \"5+67=\", \"23+9=\", \"9+13=\"),\n    @(\"27+46=\", \"7+45=\", \"74-57=\", \"32-6=\", \"88-59=\"),\n    @(\"58+28=\", \"12+69=\", \"49+42=\", \"42-39=\", \"37+46=\"),\n    @(\"54+37=\", \"36-7=\", \"27+25=\", \"75-68=\", \"8+84=\"),\n    @(\"56-28=\", \"58+34=\", \"49+35=\", \"12-7=\", \"93-69=\"),\n    @(\"76+8=\", \"71-67=\", \"27-19=\", \"29+36=\", \"74-28=\"),\n    @(\"74-6=\", \"8+56=\", \"29+35=\", \"38+15=\", \"68-49=\"),\n    @(\"44-18=\", \"98-79=\", \"36+15=\", \"70-33=\", \"80-48=\"),\n    @(\"44-19=\", \"34+37=\", \"90-65=\", \"16+58=\", \"73-15=\"),\n    @(\"86+8=\", \"83-66=\", \"25+57=\", \"37+7=\", \"9+84=\"),\n    @(\"16+19=\", \"23-6=\", \"17+69=\", \"15+78=\", \"42-26=\"),\n    @(\"46+15=\", \"55+27=\", \"79+3=\", \"52-38=\", \"71-22=\"),\n    @(\"55+39=\", \"53-35=\", \"97-29=\", \"13+58=\", \"4+17=\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
